# Update gh-pages to output generated at 456a3b4
# Applies numeric updates to column F ("浏览/热度" style counter column)
# across the "展览", "本地生活" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 5875
$ws1.Range("F15").Value = 59
$ws1.Range("F19").Value = 3962
$ws1.Range("F22").Value = 5497
$ws1.Range("F24").Value = 2169
$ws1.Range("F27").Value = 8198
$ws1.Range("F29").Value = 2223
$ws1.Range("F30").Value = 2248
$ws1.Range("F31").Value = 1349
$ws1.Range("F44").Value = 186
$ws1.Range("F46").Value = 2184

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 611
$ws3.Range("F3").Value = 792

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 611
$ws4.Range("F6").Value  = 792
$ws4.Range("F7").Value  = 5875
$ws4.Range("F14").Value = 59
$ws4.Range("F19").Value = 3962
$ws4.Range("F23").Value = 5497
$ws4.Range("F25").Value = 2169
$ws4.Range("F28").Value = 8198
$ws4.Range("F30").Value = 2223
$ws4.Range("F31").Value = 2248
$ws4.Range("F42").Value = 186
$ws4.Range("F44").Value = 2184
